$d = $word.ActiveDocument

# --- Locate the target paragraph ("{{ALLE_PROJEKTVERFASSER_NAME_ADRESSE | multiline}}") ---
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*PROJEKTVERFASSER_NAME_ADRESSE*") {
        $targetPara = $p
        break
    }
}
if ($targetPara -eq $null) {
    throw "Could not find the ALLE_PROJEKTVERFASSER_NAME_ADRESSE paragraph"
}

# --- Locate the table that immediately follows it (the EINSPRECHENDE for-loop table) ---
$targetTable = $null
foreach ($tbl in $d.Tables) {
    if ($tbl.Range.Start -eq $targetPara.Range.End) {
        $targetTable = $tbl
        break
    }
}
if ($targetTable -eq $null) {
    throw "Could not find the EINSPRECHENDE table following the target paragraph"
}

# --- Replace the paragraph + table with the two new, simplified paragraphs ---
$fullRange = $d.Range($targetPara.Range.Start, $targetTable.Range.End)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="AufzhlungVerfgung"/>
  </w:pPr>
  <w:r><w:t>{{ALLE_</w:t></w:r>
  <w:r><w:t>PROJEKTVERFASSER_NAME_ADRESSE</w:t></w:r>
  <w:r><w:t xml:space="preserve"> | </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>multiline</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>}}</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t xml:space="preserve">{% </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>for</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> POSITION in EINSPRECHENDE %}</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="AufzhlungVerfgung"/>
  </w:pPr>
  <w:r><w:t xml:space="preserve">{{POSITION.NAME}}, {{POSITION.ADRESSE}} </w:t></w:r>
  <w:r>
    <w:rPr><w:highlight w:val="yellow"/></w:rPr>
    <w:t>(inkl. Beilagen gem&#228;ss Ziffer (</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr><w:highlight w:val="yellow"/></w:rPr>
    <w:t>Nr</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr><w:highlight w:val="yellow"/></w:rPr>
    <w:t>)</w:t>
  </w:r>
  <w:r><w:t>)</w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r><w:t>{%</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>endfor</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> %}</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$fullRange.InsertXML($xml)

# --- The old (now orphaned, duplicate) table is left behind right after the
#     freshly inserted paragraphs; remove it. ---
$staleTable = $null
foreach ($tbl in $d.Tables) {
    if ($tbl.Range.Text -like "*EINSPRECHENDE*") {
        $staleTable = $tbl
        break
    }
}
if ($staleTable -ne $null) {
    $staleTable.Delete()
}

# --- The paragraph that used to follow the table ("{{GEMEINDE_NAME_ADRESSE}} ...")
#     keeps its own numbered-list formatting, but now also gets
#     w:after="0" added to its spacing (matching the paragraph that used
#     to precede the deleted table). ---
$gemeindePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*GEMEINDE_NAME_ADRESSE*") {
        $gemeindePara = $p
        break
    }
}
if ($gemeindePara -eq $null) {
    throw "Could not find the GEMEINDE_NAME_ADRESSE paragraph"
}
$gemeindePara.Format.SpaceAfter = 0
